# Auto-generated edit script applying the Ultima_Profits profitability refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2100.25
$ws.Range("I70").Value = 1200.5
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 3601.5
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -3331.5
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 2100.25
$ws.Range("I73").Value = 1200.5
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 3601.5
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -2665.5
$ws.Range("N73").Value = -10872
$ws.Range("H100").Value = 2644.9375
$ws.Range("I100").Value = 1802.375
$ws.Range("J100").Value = 3487.5
$ws.Range("K100").Value = 1802.375
$ws.Range("L100").Value = 3487.5
$ws.Range("M100").Value = -1261.375
$ws.Range("N100").Value = -4569.5
$ws.Range("H113").Value = 2743.0833
$ws.Range("I113").Value = 2487.375
$ws.Range("J113").Value = 3254.5
$ws.Range("K113").Value = 2487.375
$ws.Range("L113").Value = 3254.5
$ws.Range("M113").Value = 766.625
$ws.Range("N113").Value = -9762.5
$ws.Range("H116").Value = 3849.75
$ws.Range("J116").Value = 3999.5
$ws.Range("L116").Value = 3999.5
$ws.Range("N116").Value = -10883.5
$ws.Range("H127").Value = 899.2222
$ws.Range("I127").Value = 749.75
$ws.Range("J127").Value = 925.2174
$ws.Range("K127").Value = 2249.25
$ws.Range("L127").Value = 2775.6522
$ws.Range("M127").Value = 2710.75
$ws.Range("N127").Value = -12695.6522
$ws.Range("H132").Value = 4305.2
$ws.Range("I132").Value = 4026.25
$ws.Range("K132").Value = 12078.75
$ws.Range("M132").Value = -9548.75
$ws.Range("H133").Value = 45152
$ws.Range("J133").Value = 45152
$ws.Range("L133").Value = 45152
$ws.Range("N133").Value = -55272
$ws.Range("H137").Value = 13335281
$ws.Range("I137").Value = 1217
$ws.Range("K137").Value = 3651
$ws.Range("M137").Value = -1101
$ws.Range("H138").Value = 3207.8262
$ws.Range("I138").Value = 2675.5557
$ws.Range("J138").Value = 3550
$ws.Range("K138").Value = 8026.6671
$ws.Range("L138").Value = 10650
$ws.Range("M138").Value = -2886.6671
$ws.Range("N138").Value = -20930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6771.9307
$ws.Range("I32").Value = 7276.783
$ws.Range("K32").Value = 7276.783
$ws.Range("M32").Value = -6989.783
$ws.Range("H45").Value = 1701.9166
$ws.Range("I45").Value = 1667.1786
$ws.Range("J45").Value = 1823.5
$ws.Range("K45").Value = 1667.1786
$ws.Range("L45").Value = 1823.5
$ws.Range("M45").Value = -1290.1786
$ws.Range("N45").Value = -2577.5
$ws.Range("H63").Value = 2577.2727
$ws.Range("I63").Value = 2420
$ws.Range("J63").Value = 2996.6667
$ws.Range("K63").Value = 2420
$ws.Range("L63").Value = 2996.6667
$ws.Range("M63").Value = -1734
$ws.Range("N63").Value = -4368.6667
$ws.Range("H66").Value = 2577.2727
$ws.Range("I66").Value = 2420
$ws.Range("J66").Value = 2996.6667
$ws.Range("K66").Value = 12100
$ws.Range("L66").Value = 14983.3335
$ws.Range("M66").Value = -8668
$ws.Range("N66").Value = -21847.3335
$ws.Range("H122").Value = 5344.577
$ws.Range("I122").Value = 6522.4
$ws.Range("J122").Value = 1418.5
$ws.Range("K122").Value = 19567.2
$ws.Range("L122").Value = 4255.5
$ws.Range("M122").Value = -17117.2
$ws.Range("N122").Value = -9155.5
$ws.Range("H124").Value = 18610.545
$ws.Range("J124").Value = 18610.545
$ws.Range("L124").Value = 18610.545
$ws.Range("N124").Value = -28430.545
$ws.Range("H132").Value = 10872131
$ws.Range("I132").Value = 13890842
$ws.Range("J132").Value = 4772.4
$ws.Range("K132").Value = 41672526
$ws.Range("L132").Value = 14317.2
$ws.Range("M132").Value = -41669996
$ws.Range("N132").Value = -19377.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2633.44
$ws.Range("I134").Value = 1480.3125
$ws.Range("K134").Value = 4440.9375
$ws.Range("M134").Value = -1905.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2275
$ws.Range("I62").Value = 2293.3333
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 2293.3333
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -1669.3333
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 2275
$ws.Range("I65").Value = 2293.3333
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 11466.6665
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -8346.666499999999
$ws.Range("N65").Value = -16240
$ws.Range("H105").Value = 1779.9584
$ws.Range("I105").Value = 1016.8461
$ws.Range("J105").Value = 2681.818
$ws.Range("K105").Value = 1016.8461
$ws.Range("L105").Value = 2681.818
$ws.Range("M105").Value = 730.1539
$ws.Range("N105").Value = -6175.818
$ws.Range("H107").Value = 715.2857
$ws.Range("I107").Value = 658.5454999999999
$ws.Range("J107").Value = 923.3333
$ws.Range("K107").Value = 658.5454999999999
$ws.Range("L107").Value = 923.3333
$ws.Range("M107").Value = 1261.4545
$ws.Range("N107").Value = -4763.3333
$ws.Range("H132").Value = 3190
$ws.Range("I132").Value = 2321.4
$ws.Range("K132").Value = 6964.200000000001
$ws.Range("M132").Value = -4434.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1596
$ws.Range("I5").Value = 818
$ws.Range("J5").Value = 1985
$ws.Range("K5").Value = 2454
$ws.Range("L5").Value = 5955
$ws.Range("M5").Value = -2342
$ws.Range("N5").Value = -6179
$ws.Range("H75").Value = 2858.9167
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2858.9167
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 8576.750100000001
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -10572.7501
$ws.Range("H78").Value = 2858.9167
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2858.9167
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 25730.2503
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -35714.2503
$ws.Range("H121").Value = 943.2143
$ws.Range("I121").Value = 152
$ws.Range("J121").Value = 1382.7778
$ws.Range("K121").Value = 456
$ws.Range("L121").Value = 4148.3334
$ws.Range("M121").Value = 854
$ws.Range("N121").Value = -6768.3334
$ws.Range("H122").Value = 1063.6875
$ws.Range("I122").Value = 884.1818
$ws.Range("J122").Value = 1458.6
$ws.Range("K122").Value = 7957.6362
$ws.Range("L122").Value = 13127.4
$ws.Range("M122").Value = -5507.6362
$ws.Range("N122").Value = -18027.4
$ws.Range("H135").Value = 1596
$ws.Range("I135").Value = 818
$ws.Range("J135").Value = 1985
$ws.Range("K135").Value = 7362
$ws.Range("L135").Value = 17865
$ws.Range("M135").Value = -4827
$ws.Range("N135").Value = -22935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16777746
$ws.Range("J80").Value = 3249734.8
$ws.Range("L80").Value = 3249734.8
$ws.Range("N80").Value = -3251730.8
$ws.Range("H83").Value = 16777746
$ws.Range("J83").Value = 3249734.8
$ws.Range("L83").Value = 16248674
$ws.Range("N83").Value = -16258658
$ws.Range("H97").Value = 929.4838999999999
$ws.Range("J97").Value = 686.5
$ws.Range("L97").Value = 686.5
$ws.Range("N97").Value = -1678.5
$ws.Range("H107").Value = 1757.9286
$ws.Range("I107").Value = 2086
$ws.Range("J107").Value = 937.75
$ws.Range("K107").Value = 2086
$ws.Range("L107").Value = 937.75
$ws.Range("M107").Value = -166
$ws.Range("N107").Value = -4777.75
$ws.Range("H122").Value = 4763369.5
$ws.Range("I122").Value = 5556765
$ws.Range("J122").Value = 2997.5
$ws.Range("K122").Value = 16670295
$ws.Range("L122").Value = 8992.5
$ws.Range("M122").Value = -16667845
$ws.Range("N122").Value = -13892.5
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -54900
$ws.Range("H126").Value = 5246.125
$ws.Range("I126").Value = 3665
$ws.Range("J126").Value = 5773.1665
$ws.Range("K126").Value = 10995
$ws.Range("L126").Value = 17319.4995
$ws.Range("M126").Value = -8525
$ws.Range("N126").Value = -22259.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1415.25
$ws.Range("I93").Value = 1275.6666
$ws.Range("J93").Value = 1624.625
$ws.Range("K93").Value = 1275.6666
$ws.Range("L93").Value = 1624.625
$ws.Range("M93").Value = -27.66660000000002
$ws.Range("N93").Value = -4120.625
$ws.Range("H139").Value = 44335.4
$ws.Range("J139").Value = 44744.89
$ws.Range("L139").Value = 44744.89
$ws.Range("N139").Value = -55024.89

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1246.614
$ws.Range("I132").Value = 936.88635
$ws.Range("J132").Value = 2294.923
$ws.Range("K132").Value = 2810.65905
$ws.Range("L132").Value = 6884.768999999999
$ws.Range("M132").Value = -280.6590500000002
$ws.Range("N132").Value = -11944.769
